$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("discription") on each data row is repurposed from a duplicate
# of the title into a short "how to play" (玩法) instruction string.
$ws.Range("C2").Value = "玩法: 上載一張自己的相，就可轉為吉卜力(Ghibli)風格"
$ws.Range("C3").Value = "玩法: 上載一張成年, 一張小時候的相，就可得出自己的合照"
$ws.Range("C4").Value = "玩法: 上載一張自己的相, 可生成自己的卡通貼紙"
$ws.Range("C5").Value = "玩法: 上載一張自己的相和一張衣服相，就可更換衣服"
$ws.Range("C6").Value = "玩法: 上載一張自己的相, 就可做出一副3D figurine圖"
$ws.Range("C7").Value = "玩法: 上載一張自己的相, 和一枝國旗的相，就可向國家致敬"
$ws.Range("C8").Value = "玩法: 上載一張有天空的相，大大的月亮就可出現"
$ws.Range("C9").Value = "玩法: 上載一張自己的相, ，就可以與暴龍合照"

# Row 5 ("更換衣服" / change clothes) gains a second reference-image cell
# (F5) naming the clothing photo, alongside the existing self-photo (E5).
# Copy E5's formatting onto F5 so the new cell carries the same cell style
# (s="1") as every other populated cell, rather than the blank default.
$ws.Range("F5").Value = "dress.jpg"
$ws.Range("E5").Copy() | Out-Null
$ws.Range("F5").PasteSpecial(-4122) | Out-Null
